$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing F column paths (rows 2-6) and add E column scores
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "C:\Users\GreaTiger\Desktop\project\results\السرية الثامنة\1_result.jpg"

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = "C:\Users\GreaTiger\Desktop\project\results\السرية الثامنة\2_result.jpg"

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "C:\Users\GreaTiger\Desktop\project\results\السرية الثامنة\3_result.jpg"

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "C:\Users\GreaTiger\Desktop\project\results\السرية الثامنة\4_result.jpg"

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = "C:\Users\GreaTiger\Desktop\project\results\السرية الثامنة\5_result.jpg"

# New rows 7-8: A/B are stored as text (matching the existing A2:B6 convention),
# so force text format before assigning the numeric-looking strings.
$ws.Range("A7:B8").NumberFormat = "@"

# Add new row 7
$ws.Range("A7").Value = "6"
$ws.Range("B7").Value = "490"
$ws.Range("C7").Value = "نورالدين عبدالحميد"
$ws.Range("D7").Value = "طالب"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = "C:\Users\GreaTiger\Desktop\project\results\السرية الثامنة\6_result.jpg"

# Add new row 8
$ws.Range("A8").Value = "7"
$ws.Range("B8").Value = "505"
$ws.Range("C8").Value = "محمد صبحى احمد"
$ws.Range("D8").Value = "طالب"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = "C:\Users\GreaTiger\Desktop\project\results\السرية الثامنة\7_result.jpg"
